$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Day" column (B) for every 2022 row (B2195:B2558) was stuck at 1.
# Data has now been updated through 1/6/2022, so the day-of-year counter
# needs to increment properly: B2194 (1/1/2022) stays 1, and each
# following row counts up by one (B2195 = 2, B2196 = 3, ..., B2558 = 365).
for ($r = 2195; $r -le 2558; $r++) {
    $ws.Cells.Item($r, 2).Value = $r - 2193
}

# Leave the sheet scrolled to / selected on the range that was just
# corrected, matching where the editor's cursor ended up.
$ws.Range("B2194:B2558").Select()
